# Add a new "MS_DEF" column (F) to the mapping sheet, filling every
# existing data row (2-27) with the default value "[]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the style already used by the other headers
# (B1:E1) - copy that formatting onto F1 rather than re-building it by hand.
$headerCell = $ws.Range("F1")
$headerCell.Value = "MS_DEF"
$ws.Range("E1").Copy() | Out-Null
$headerCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new column for every data row with the empty-list placeholder.
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 6).Value = "[]"
}
